$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.770.00"
$ws.Range("E2").Value = "  -0.56%  "

# Row 3
$ws.Range("D3").Value = "3.507.59"
$ws.Range("E3").Value = "  -1.13%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.67"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -1.71%  "

# Row 6
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.75"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -0.60%  "

# Row 7
$ws.Range("D7").Value = "3.508.66"
$ws.Range("E7").Value = "  -1.09%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("E9").Value = "  -0.83%  "

# Row 10
$ws.Range("E10").Value = "  +1.69%  "

# Row 11
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.12"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.44%  "

# Row 12
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("D13").Value = "4.104.14"
$ws.Range("E13").Value = "  -1.15%  "

# Row 14
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.73"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +3.65%  "

# Row 15
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -1.02%  "

# Row 16
$ws.Range("E16").Value = "  +0.79%  "

# Row 17
$ws.Range("D17").Value = "3.507.77"
$ws.Range("E17").Value = "  -1.25%  "

# Row 18
$ws.Range("D18").Value = "64.779.87"
$ws.Range("E18").Value = "  -0.78%  "

# Row 19
$ws.Range("E19").Value = "  +0.67%  "

# Row 20
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.25"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -0.67%  "

# Row 21
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.68"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -2.44%  "

# Row 22
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.54"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +0.32%  "

# Row 23
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.576"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").Value = "3.648.83"
$ws.Range("E24").Value = "  -1.24%  "

# Row 25
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.26"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +0.42%  "

# Row 26
$ws.Range("E26").Value = "  -0.09%  "

# Row 27
$ws.Range("E27").Value = "  -2.88%  "

# Row 28
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.58"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +3.40%  "

# Row 29
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +0.04%  "

# Row 30
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.44"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -3.50%  "

# Row 31
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.26"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -0.87%  "

# Row 32
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.22"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -3.19%  "

# Row 33
$ws.Range("D33").Value = "3.512.81"
$ws.Range("E33").Value = "  -1.20%  "

# Row 34
$ws.Range("E34").Value = "  +0.00%  "

# Row 35
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.97"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -0.29%  "

# Row 36
$ws.Range("E36").Value = "  -0.54%  "

# Row 37
$ws.Range("E37").Value = "  +4.04%  "

# Row 38
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "171.89"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +1.15%  "

# Row 39
$ws.Range("E39").Value = "  +4.04%  "

# Row 40
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.99"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +1.38%  "

# Row 41
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0808"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -0.08%  "

# Row 42
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.815"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  -1.32%  "

# Row 43
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.42"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +0.92%  "

# Row 44
$ws.Range("E44").Value = "  -0.09%  "

# Row 45
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.35"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -1.63%  "

# Row 46
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.23"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -0.77%  "

# Row 47
$ws.Range("E47").Value = "  -0.34%  "

# Row 48
$ws.Range("E48").Value = "  -0.10%  "

# Row 49
$ws.Range("D49").Value = "2.482.87"
$ws.Range("E49").Value = "  +1.05%  "

# Row 50
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.87"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -0.36%  "

# Row 51
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.904"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +4.03%  "
